# Append the new draw result as row 94 on the "Results" sheet.
# Source data (same shape/format as every prior row):
#   A: Date            2025-12-19
#   B: Game             Pick 4
#   C: Phase            251219
#   D: Result           8-3-3-8
#   E: InsertedAt       2025-12-19T21:37:22.987+04:00
#
# Every existing cell in the sheet is stored as plain TEXT (t="str"/shared
# string), even values that look like dates or numbers (e.g. "251219").
# A straight `.Value = "..."` assignment lets Excel's type inference kick
# in and would silently turn "2025-12-19" into a date serial and "251219"
# into a number (each picking up a new number-format style along the way).
# To keep the new row textual - matching the rest of the column - and to
# avoid minting any new cell styles, write the "risky" cells as formulas
# that evaluate to text, then convert the whole row to static values via
# Copy / Paste-Special-Values (which preserves the text result without
# re-triggering Excel's "does this look like a number/date" auto-convert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 94

$ws.Cells.Item($row, 1).Formula = '="2025-12-19"'
$ws.Cells.Item($row, 2).Value = "Pick 4"
$ws.Cells.Item($row, 3).Formula = '="251219"'
$ws.Cells.Item($row, 4).Value = "8-3-3-8"
$ws.Cells.Item($row, 5).Value = "2025-12-19T21:37:22.987+04:00"

$newRowRange = $ws.Range("A" + $row + ":E" + $row)
$newRowRange.Copy()
$newRowRange.PasteSpecial(-4163)  # xlPasteValues - freeze formulas to static text values
